$d = $word.ActiveDocument

# Change 1: merge " que" / " lea, dende o teclado, " runs and
# drop the gramStart/gramEnd proofErr markers around them.
$d.Content.Find.Execute(
    " que lea, dende o teclado, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " que lea, dende o teclado, ",
    2)

# Change 2: merge "desglose en billetes de 100, 20, 5 e " / "moedas" / " de 1 € ."
# runs and drop the spellStart/spellEnd proofErr markers around "moedas".
$d.Content.Find.Execute(
    "desglose en billetes de 100, 20, 5 e moedas de 1 € .",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "desglose en billetes de 100, 20, 5 e moedas de 1 € .",
    2)
